# Update the "Förändrad" (changed) date in column C for every data row,
# and append the row's "Beteckning" (column A) as the friendly-name second
# argument to every HYPERLINK() formula (columns S, T, U, V, W, X, Y).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 358
$linkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {
    # Column C: "Förändrad" date serial 45184 -> 45186
    $ws.Range("C$r").Value = 45186

    # Friendly text used as the HYPERLINK() second argument is the
    # row's "Beteckning" value in column A.
    $label = $ws.Range("A$r").Text

    foreach ($col in $linkCols) {
        $cell = $ws.Range("$col$r")
        if ($cell.HasFormula) {
            $f = $cell.Formula
            $hasLabelAlready = $f.Contains(',"' + $label + '")') -or $f.Contains(', "' + $label + '")')
            if ($f.StartsWith('=HYPERLINK(') -and $f.EndsWith('")') -and -not $hasLabelAlready) {
                # Only add the friendly-name argument if it isn't already there.
                $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $label + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
